$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix misspelled Russian label ("Сирийный номер" -> "Серийный номер", i.e.
# "Syrian number" typo corrected to "Serial number") on the connection
# settings sheet, as part of tightening up error checking around it.
$ws.Range("H5").Value = "Серийный номер"

# Move the active selection to the corrected cell.
$ws.Range("H5").Select()
